$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the new "Escalabilidad" (Scalability) non-functional requirement row
$ws.Range("B9").Value = "Cantidad de usuarios"
$ws.Range("C9").Value = "Contemplar que el sistema soporte la cantidad de usuarios necesarios sin afectar el desempeño de la aplcacion"
$ws.Range("D9").Value = "Escalabilidad"

# Match the wrap-text style used by the other description cells in column C
$ws.Range("C9").WrapText = $true

# Adjust row height to fit the new wrapped text, matching the saved workbook
$ws.Rows.Item(9).RowHeight = 30.75

# Update the view state to reflect where the user was working
$ws.Range("B9").Select()
